# Auto-generated script: update the 25 data cells of the single table
# in the document by replacing each cell's text content in place,
# preserving run/paragraph formatting (rFonts, sz, jc, etc.).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "14÷8=1, 6"
$t.Cell(1, 2).Range.Text = "93÷3=31, 0"
$t.Cell(1, 3).Range.Text = "25÷2=12, 1"
$t.Cell(1, 4).Range.Text = "75÷8=9, 3"
$t.Cell(1, 5).Range.Text = "89÷8=11, 1"

# Row 5
$t.Cell(5, 1).Range.Text = "42÷7=6, 0"
$t.Cell(5, 2).Range.Text = "47÷7=6, 5"
$t.Cell(5, 3).Range.Text = "17÷4=4, 1"
$t.Cell(5, 4).Range.Text = "38÷4=9, 2"
$t.Cell(5, 5).Range.Text = "59÷8=7, 3"

# Row 9
$t.Cell(9, 1).Range.Text = "54÷8=6, 6"
$t.Cell(9, 2).Range.Text = "81÷6=13, 3"
$t.Cell(9, 3).Range.Text = "13÷4=3, 1"
$t.Cell(9, 4).Range.Text = "18÷7=2, 4"
$t.Cell(9, 5).Range.Text = "15÷4=3, 3"

# Row 13
$t.Cell(13, 1).Range.Text = "51÷6=8, 3"
$t.Cell(13, 2).Range.Text = "98÷5=19, 3"
$t.Cell(13, 3).Range.Text = "17÷3=5, 2"
$t.Cell(13, 4).Range.Text = "40÷5=8, 0"
$t.Cell(13, 5).Range.Text = "92÷7=13, 1"

# Row 17
$t.Cell(17, 1).Range.Text = "24÷3=8, 0"
$t.Cell(17, 2).Range.Text = "52÷2=26, 0"
$t.Cell(17, 3).Range.Text = "17÷5=3, 2"
$t.Cell(17, 4).Range.Text = "63÷8=7, 7"
$t.Cell(17, 5).Range.Text = "22÷8=2, 6"

